$wb = $excel.ActiveWorkbook

# Overview sheet: status text changed everywhere "Ready for handoff" appeared
$ov = $wb.Worksheets.Item("Overview")
$ov.Range("B2").Value = "Handed back: in sync with en-US"
$ov.Range("C2").Value = "Handed back: in sync with en-US"
$ov.Range("B3").Value = "Handed back: in sync with en-US"
$ov.Range("C3").Value = "Handed back: in sync with en-US"

# zh-cn sheet
$zh = $wb.Worksheets.Item("zh-cn")
$zh.Range("B2").Value = "Handed back: in sync with en-US"
$zh.Range("B3").Value = "Handed back: in sync with en-US"
$zh.Range("G2").Value = "2016-03-08 10:42:59"
$zh.Range("G3").Value = "2016-03-08 10:42:59"

$zh.Hyperlinks.Add($zh.Range("E2"), "https://github.com/OpenLocalizationTest/oltest/blob/04e13ba9491c405a4dd4a5fd6421fdfe320c02cc/e2e/3ca43955-bd38-4955-aafa-69ccbe78ff67.md", "", "", "3ca43955-bd38-4955-aafa-69ccbe78ff67.md")
$zh.Hyperlinks.Add($zh.Range("F2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/0dcfe663ff287638f1e7daf8903e2c9222e7438a/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/3ca43955-bd38-4955-aafa-69ccbe78ff67.8b726d99015d6d04615f5d4a555d51fa6ca19b07.zh-cn.xlf", "", "", "3ca43955-bd38-4955-aafa-69ccbe78ff67.8b726d99015d6d04615f5d4a555d51fa6ca19b07.zh-cn.xlf")
$zh.Hyperlinks.Add($zh.Range("E3"), "https://github.com/OpenLocalizationTest/oltest/blob/04e13ba9491c405a4dd4a5fd6421fdfe320c02cc/e2e/bbeb3347-5d2c-450a-b41a-f9167885289b.md", "", "", "bbeb3347-5d2c-450a-b41a-f9167885289b.md")
$zh.Hyperlinks.Add($zh.Range("F3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/0dcfe663ff287638f1e7daf8903e2c9222e7438a/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/bbeb3347-5d2c-450a-b41a-f9167885289b.1a9e9f6d8ea9306613957ff624dc3b666edb2bf1.zh-cn.xlf", "", "", "bbeb3347-5d2c-450a-b41a-f9167885289b.1a9e9f6d8ea9306613957ff624dc3b666edb2bf1.zh-cn.xlf")

# de-de sheet
$de = $wb.Worksheets.Item("de-de")
$de.Range("B2").Value = "Handed back: in sync with en-US"
$de.Range("B3").Value = "Handed back: in sync with en-US"
$de.Range("G2").Value = "2016-03-08 10:43:09"
$de.Range("G3").Value = "2016-03-08 10:43:09"

$de.Hyperlinks.Add($de.Range("E2"), "https://github.com/OpenLocalizationTest/oltest/blob/04e13ba9491c405a4dd4a5fd6421fdfe320c02cc/e2e/3ca43955-bd38-4955-aafa-69ccbe78ff67.md", "", "", "3ca43955-bd38-4955-aafa-69ccbe78ff67.md")
$de.Hyperlinks.Add($de.Range("F2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/39c69f6d5f07303fc309ba23d00ae24b55eacf66/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/3ca43955-bd38-4955-aafa-69ccbe78ff67.8b726d99015d6d04615f5d4a555d51fa6ca19b07.de-de.xlf", "", "", "3ca43955-bd38-4955-aafa-69ccbe78ff67.8b726d99015d6d04615f5d4a555d51fa6ca19b07.de-de.xlf")
$de.Hyperlinks.Add($de.Range("E3"), "https://github.com/OpenLocalizationTest/oltest/blob/04e13ba9491c405a4dd4a5fd6421fdfe320c02cc/e2e/bbeb3347-5d2c-450a-b41a-f9167885289b.md", "", "", "bbeb3347-5d2c-450a-b41a-f9167885289b.md")
$de.Hyperlinks.Add($de.Range("F3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/39c69f6d5f07303fc309ba23d00ae24b55eacf66/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/bbeb3347-5d2c-450a-b41a-f9167885289b.1a9e9f6d8ea9306613957ff624dc3b666edb2bf1.de-de.xlf", "", "", "bbeb3347-5d2c-450a-b41a-f9167885289b.1a9e9f6d8ea9306613957ff624dc3b666edb2bf1.de-de.xlf")
